# Update task-order worksheets: rename sheets (new timestamps) and
# refresh the stimulus-file / condition values logged in column B,
# improving the accuracy of stimulus presentation time-logging.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512554812004519"
$ws1.Range("B2").Value = "go_stims-1651255481167254.csv"
$ws1.Range("B3").Value = "GNG_stims-16512554811833642.csv"
$ws1.Range("B4").Value = "go_stims-16512554811852875.csv"
$ws1.Range("B5").Value = "GNG_stims-16512554811994097.csv"

# --- Sheet 2: NB -----------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16512554828371131"
$ws2.Range("B2").Value = "TB-1651255482814433.csv"
$ws2.Range("B3").Value = "ZB-match_5-16512554818761458.csv"
$ws2.Range("B4").Value = "OB-16512554822250252.csv"
$ws2.Range("B5").Value = "TB-16512554827413402.csv"
$ws2.Range("B6").Value = "ZB-match_2-1651255481424914.csv"
$ws2.Range("B7").Value = "OB-16512554823695858.csv"
$ws2.Range("B8").Value = "TB-16512554825074532.csv"
$ws2.Range("B9").Value = "OB-1651255482003001.csv"
$ws2.Range("B10").Value = "ZB-match_9-16512554816314726.csv"

# --- Sheet 3: RS -------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651255482841217"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16512554829108694"
$ws4.Range("B2").Value = "MM_stims-16512554828659723.csv"
$ws4.Range("B3").Value = "ZM_stims-16512554828451188.csv"
$ws4.Range("B4").Value = "MM_stims-16512554828976758.csv"
$ws4.Range("B5").Value = "ZM_stims-16512554828679621.csv"
$ws4.Range("B6").Value = "MM_stims-16512554829098735.csv"
$ws4.Range("B7").Value = "ZM_stims-16512554828986802.csv"

# --- Sheet 5: vSAT -------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1651255483017589"
$ws5.Range("B2").Value = "vSAT_stims-1651255483002457.csv"
$ws5.Range("B3").Value = "SAT_stims-16512554829178002.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512554829569082.csv"
$ws5.Range("B5").Value = "SAT_stims-16512554829343317.csv"
